# Doing Updates for Financials
# Refresh the yearly financial data pulled from the data source: most rows
# shift one column to the left (oldest year drops off, and the newest
# column either gets a freshly-pulled number or "NA" when no figure is
# available yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KIRK")

# Earnings Before Interest And Taxes (row 21)
$ws.Range("D21").Value = 35400
$ws.Range("E21").Value = 39400
$ws.Range("F21").Value = 45100
$ws.Range("G21").Value = 45000
$ws.Range("H21").Value = 37400
$ws.Range("I21").Value = 34200
$ws.Range("J21").Value = "NA"

# Depreciation (row 83)
$ws.Range("D83").Value = 25300
$ws.Range("E83").Value = 22200
$ws.Range("F83").Value = 18600
$ws.Range("G83").Value = 15900
$ws.Range("H83").Value = 13200
$ws.Range("I83").Value = 12400
$ws.Range("J83").Value = "NA"

# Total Cash Flow From Operating Activities (row 89)
$ws.Range("D89").Value = 51900
$ws.Range("E89").Value = 33200
$ws.Range("F89").Value = 44500
$ws.Range("G89").Value = 39200
$ws.Range("H89").Value = 32300
$ws.Range("I89").Value = 41800
$ws.Range("J89").Value = 36700

# Capital Expenditures (row 91)
$ws.Range("D91").Value = -32200
$ws.Range("E91").Value = -35100
$ws.Range("F91").Value = -29600
$ws.Range("G91").Value = -18000
$ws.Range("H91").Value = -31400
$ws.Range("I91").Value = -26700
$ws.Range("J91").Value = -22600

# Total Cash Flows From Investing Activities (row 94)
$ws.Range("D94").Value = -32200
$ws.Range("E94").Value = -35100
$ws.Range("F94").Value = -29600
$ws.Range("G94").Value = -18000
$ws.Range("H94").Value = -31400
$ws.Range("I94").Value = -26700
$ws.Range("J94").Value = "NA"

# Dividends Paid (row 96)
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = -26000
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0

# Other Cash Flows from Financing Activities (row 100)
$ws.Range("D100").Value = -200
$ws.Range("E100").Value = -52800
$ws.Range("F100").Value = -4800
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = -16300
$ws.Range("I100").Value = -23200
$ws.Range("J100").Value = "NA"

# Total Cash Flows From Financing Activities (row 101)
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = "NA"

# Effect Of Exchange Rate Changes (row 102)
$ws.Range("D102").Value = 19600
$ws.Range("E102").Value = -54800
$ws.Range("F102").Value = 10100
$ws.Range("G102").Value = 21300
$ws.Range("H102").Value = -15300
$ws.Range("I102").Value = -8100
$ws.Range("J102").Value = 14800
